$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 14) down to the two
# new rows (15 and 16) so they pick up the same styles / borders as Excel
# would apply when the table grows.
$ws.Range("A14:H14").Copy()
$ws.Range("A15:H16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 15
$ws.Range("A15").Value = 45057
$ws.Range("B15").Value = 0.56527777777777777
$ws.Range("C15").Value = 0.65208333333333335
$ws.Range("D15").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E15").Value = "CPNV"
$ws.Range("F15").Value = "Base de données"
$ws.Range("G15").Value = "M. Meylan"
$ws.Range("H15").Value = "J'ai fini le scriptqui crée la base de données"

# Row 16
$ws.Range("A16").Value = 45057
$ws.Range("B16").Value = 0.65277777777777779
$ws.Range("C16").Value = 0.66041666666666665
$ws.Range("D16").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E16").Value = "CPNV"
$ws.Range("F16").Value = "Base de données"
$ws.Range("G16").Value = "J'ai fait un script qui crée un backup de la base de données"

# Grow the table (ListObject) and its AutoFilter range to include the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H16"))

# Update the active selection like Excel would after editing the last cell.
$ws.Range("H16").Select() | Out-Null
